$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 40: AVDD/VDD note with reference link ---
$ws.Range("A40").Value = "AVDD VDD用100uH磁珠連接"
$ws.Range("D40").Value = "https://blog.csdn.net/D_Katter/article/details/127743365"

# --- Ver0. problem section: mark items as solved / add "V" column ---
$ws.Range("B45").Value = "solved"

$ws.Range("B46").Value = "V"
$ws.Range("B47").Value = "V"
$ws.Range("B48").Value = "V"
$ws.Range("B49").Value = "V"
$ws.Range("B50").Value = "V"

# Rows 52/53 (before) become rows 51/52 (after); clear the old pair then
# rewrite content one row higher, still flagged solved ("V").
$oldA52 = $ws.Range("A52").Value()
$oldA53 = $ws.Range("A53").Value()

$ws.Range("A52").ClearContents()
$ws.Range("B52").ClearContents()
$ws.Range("A53").ClearContents()
$ws.Range("B53").ClearContents()

$ws.Range("A51").Value = $oldA52
$ws.Range("B51").Value = "V"
$ws.Range("A52").Value = $oldA53
$ws.Range("B52").Value = "V"

# --- New rows 55-56: further solved issues ---
$ws.Range("A55").Value = "37pin D-sub確認固定孔位置"
$ws.Range("B55").Value = "V"

$ws.Range("A56").Value = "top paste 有點不對, layout用藍色那層的不會裸露出來"
$ws.Range("B56").Value = "V"

# --- New section header row 60: "final check" (reuse bold/yellow style from A45) ---
$ws.Range("A45").Copy()
$ws.Range("A60").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A60").Value = "final check"

$ws.Range("A61").Value = "solder PAD and hole PAD show up"
$ws.Range("B61").Value = "V"

$ws.Range("A62").Value = "Ground, small broken island check after pour polygon"
$ws.Range("B62").Value = "V"

$ws.Range("A63").Value = "top overlay check"
$ws.Range("B63").Value = "V"

$ws.Range("A64").Value = "top/bottom paste, top/bottom solder check"
$ws.Range("B64").Value = "V"

# --- Column A grew wider to fit the new longer text ---
$ws.Columns("A:A").ColumnWidth = 54.25

# --- Update view: scroll down a bit and select C64 ---
$ws.Range("C64").Select()
